$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "30.466.52"
$ws.Cells.Item(2,5).Value = "  -0.78%  "

$ws.Cells.Item(3,4).Value = "2.088.84"
$ws.Cells.Item(3,5).Value = "  -1.27%  "

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "1.003"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = "  +0.56%  "

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "329.11"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "  -0.40%  "

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "1.002"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "  +0.61%  "

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "0.5213"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = "  +0.89%  "

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.4345"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "  -0.65%  "

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "51.94"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "  +13.57%  "

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.08834"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "  -2.46%  "

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "1.155"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "  -1.89%  "

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "24.47"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "  -2.53%  "

$ws.Cells.Item(13,4).Value = "2.103.12"
$ws.Cells.Item(13,5).Value = "  -0.18%  "

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "6.690"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "  -1.38%  "

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "7.676"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "  +0.76%  "

$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "95.66"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = "  -1.83%  "

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "1.004"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = "  +0.71%  "

$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "0.00001119"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "  -1.77%  "

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "0.06609"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "  +0.23%  "

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "19.22"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "  +0.10%  "

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "1.002"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "  +0.65%  "

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "6.280"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "  -1.91%  "

$ws.Cells.Item(23,4).Value = "30.513.76"
$ws.Cells.Item(23,5).Value = "  -1.55%  "

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "12.21"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = "  +1.83%  "

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "2.341"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "  +3.08%  "

$ws.Cells.Item(26,4).Value = "2.335.27"
$ws.Cells.Item(26,5).Value = "  -0.72%  "

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "22.26"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "  -2.76%  "

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "2.589"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = "  +0.65%  "

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "162.06"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = "  -0.90%  "

$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "131.75"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = "  -1.45%  "

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "1.199"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = "  -0.35%  "

$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "0.1066"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = "  +0.00%  "

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "1.661"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = "  +6.82%  "

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "6.152"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = "  -1.56%  "

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "3.893"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = "  +1.42%  "

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "10.09"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = "  +6.36%  "

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = "0.02585"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = "  +0.04%  "

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "0.06851"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = "  +0.82%  "

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "5.454"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = "  -2.87%  "

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "12.61"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "  -0.01%  "

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "0.2260"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "  +0.68%  "

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.6913"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = "  +2.07%  "

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "1.264"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = "  +1.32%  "

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "1.002"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "  +0.71%  "

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "0.6371"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "  +0.99%  "

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "13.95"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = "  -2.61%  "

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "2.205"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = "  -1.80%  "

$ws.Cells.Item(48,2).Value = "PancakeSwap"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "3.628"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = "  -0.75%  "

$ws.Cells.Item(49,2).Value = "WEMIXTOKEN"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "1.240"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = "  +14.43%  "

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "1.239"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "  -2.47%  "

$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "81.65"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "  -2.19%  "
